# Insert a new data row before row 81 (pushing existing rows 81-133 down to 82-134)
# and populate the newly inserted row 81 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81; existing row 81 (and below) shift down to 82.
$ws.Rows.Item(81).Insert()

# Copy the date number format used by the other "Fecha" cells (column D) onto the new cell.
$ws.Cells.Item(81, 4).NumberFormat = $ws.Cells.Item(82, 4).NumberFormat

# Populate the new row 81 with the record's values.
$ws.Cells.Item(81, 1).Value = 2
$ws.Cells.Item(81, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(81, 3).Value = "Coquimbo"
$ws.Cells.Item(81, 4).Value = 44762
$ws.Cells.Item(81, 5).Value = 4
$ws.Cells.Item(81, 6).Value = 100112024
$ws.Cells.Item(81, 7).Value = "Choclo"
$ws.Cells.Item(81, 8).Value = "Dulce o Americano"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 1100
$ws.Cells.Item(81, 11).Value = 25000
$ws.Cells.Item(81, 12).Value = 27000
$ws.Cells.Item(81, 13).Value = 26000
$ws.Cells.Item(81, 14).Value = "`$/malla 60 unidades"
$ws.Cells.Item(81, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(81, 16).Value = 433
$ws.Cells.Item(81, 17).Value = 60
$ws.Cells.Item(81, 18).Value = "Hortaliza"
